$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: SCD0259 -> SCD0016
$ws.Name = "SCD0016"

# Update the TC_ID cell (B2) from "DGS-274" to "SCD0016-033"
$ws.Range("B2").Value = "SCD0016-033"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns("B").ColumnWidth = 12.5703125

# Header row (row 1) becomes left+vertically centered
$headerRow = $ws.Range("A1:O1")
$headerRow.HorizontalAlignment = -4131
$headerRow.VerticalAlignment = -4108

# B2's TC_ID now uses the bigger Arial font (size 10) used elsewhere, left/center aligned
$b2 = $ws.Range("B2")
$b2.Font.Name = "Arial"
$b2.Font.Size = 10
$b2.HorizontalAlignment = -4131
$b2.VerticalAlignment = -4108

# Left-align the rest of the data row (row 2) cells that switched from default/right to left
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("C2:E2").HorizontalAlignment = -4131
$ws.Range("F2").HorizontalAlignment = -4131
$ws.Range("G2").HorizontalAlignment = -4131
$ws.Range("H2:L2").HorizontalAlignment = -4131
$ws.Range("M2").HorizontalAlignment = -4131
$ws.Range("O2").HorizontalAlignment = -4131

# N2 gains explicit center/left alignment (previously had none)
$n2 = $ws.Range("N2")
$n2.HorizontalAlignment = -4131
$n2.VerticalAlignment = -4108

# Adjust the window view: no frozen/scrolled top-left cell, selection moves to B3
$excel.ActiveWindow.TopLeftCell = $ws.Range("A1")
$ws.Range("B3").Select()
